# Saldo_guide.xlsx update
# 1. Rename the worksheet (was the raw export name, now "Saldo_guide")
# 2. Bump the "Dt. Referencia" (column G) date from 45428 -> 45429 for every data row
# 3. A handful of rows also got their Saldo Previsto / Vl. Projetado / Vl. Total
#    (D / E / H) figures corrected

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Sheet name -------------------------------------------------------
$ws.Name = "Saldo_guide"

# --- 2. Column G date roll-forward for all data rows (2..257) ------------
for ($r = 2; $r -le 257; $r++) {
    $ws.Cells.Item($r, 7).Value = 45429
}

# --- 3. Row-specific corrections to D / E / H -----------------------------
# Row 5
$ws.Range("D5").Value = 8148.72
$ws.Range("E5").Value = 0
$ws.Range("H5").Value = 8148.72

# Row 42
$ws.Range("D42").Value = 513.73
$ws.Range("H42").Value = 513.73

# Row 51
$ws.Range("D51").Value = 29957.62
$ws.Range("H51").Value = 29957.62

# Row 54
$ws.Range("D54").Value = 119.34
$ws.Range("H54").Value = 119.34

# Row 57
$ws.Range("D57").Value = 189.89
$ws.Range("H57").Value = 189.89

# Row 96
$ws.Range("D96").Value = 0
$ws.Range("H96").Value = 0

# Row 98
$ws.Range("D98").Value = 99.65
$ws.Range("E98").Value = 0
$ws.Range("H98").Value = 99.65

# Row 103
$ws.Range("D103").Value = -0.24
$ws.Range("E103").Value = 0
$ws.Range("H103").Value = -0.24

# Row 226
$ws.Range("D226").Value = 99.68
$ws.Range("E226").Value = 0
$ws.Range("H226").Value = 99.68
